$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("E2").Value = 25.74000000000058
$ws.Range("H2").Value = [double]"2.046494054608583e-16"
$ws.Range("K2").Value = 46.33824671219084
$ws.Range("L2").Value = "[41.37069612838211, 51.30579729599957]"
$ws.Range("O2").Value = 1.352237078121732
$ws.Range("P2").Value = "[1.2390265320464247, 1.4654476241970391]"
$ws.Range("S2").Value = 54.47010478739367
$ws.Range("T2").Value = "[51.45938144078063, 57.480828134006714]"
$ws.Range("W2").Value = 20.20036036036082
$ws.Range("X2").Value = 19.73657657657703
$ws.Range("Y2").Value = 20.66414414414461

# Row 3 changes
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 22
$ws.Range("H3").Value = [double]"2.046494054608583e-16"
$ws.Range("I3").Value = 0.9572681818530997
$ws.Range("K3").Value = 48.55320690959655
$ws.Range("L3").Value = "[40.56630092373977, 56.54011289545334]"
$ws.Range("O3").Value = 0.1446579199851161
$ws.Range("P3").Value = "[-0.03144737390980623, 0.3207632138800385]"
$ws.Range("Q3").Value = 0.1066289393885655
$ws.Range("R3").Value = 0.1066289393885655
$ws.Range("S3").Value = 54.68317344460878
$ws.Range("T3").Value = "[49.69041155217134, 59.675935337046226]"
$ws.Range("W3").Value = 21.49349349349349
$ws.Range("X3").Value = 20.87687687687688
$ws.Range("Y3").Value = 22.11011011011011
